$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cookbook")

# The "crockpot" column (E) used to hold "Yes"/"No" text. Only rows that are
# actually crockpot recipes keep a value now ("Crockpot"); the others are
# cleared out entirely.
$ws.Range("E2").Value = "Crockpot"
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").Value = "Crockpot"
$ws.Range("E7").Value = "Crockpot"

# Insert a new "img_source" column right after "soure_type" (before the old
# "recipe"/"directions" columns, which shift right from H/I to I/J).
$ws.Columns.Item(8).Insert()
$ws.Range("H1").Value = "img_source"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = "static/example.jpg"
}
$ws.Columns.Item(8).ColumnWidth = 15.27

# Print setup for the table.
$ws.PageSetup.Orientation = 1

# Reflect the last edit (filling the new column) in the sheet's selection.
$ws.Range("H2:H7").Select()
